# Auto-generated edit script applying diff changes to cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.252.15"
$ws.Range("E2").Value = "'  -2.39%  "
$ws.Range("D3").Value = "'1.559.28"
$ws.Range("E3").Value = "'  -3.73%  "
$ws.Range("E4").Value = "'  -0.23%  "
$ws.Range("D5").Value = "'206.16"
$ws.Range("E5").Value = "'  -3.23%  "
$ws.Range("E7").Value = "'  -4.47%  "
$ws.Range("D8").Value = "'0.0613"
$ws.Range("E8").Value = "'  -0.48%  "
$ws.Range("D9").Value = "'0.242"
$ws.Range("E9").Value = "'  -2.81%  "
$ws.Range("D10").Value = "'17.73"
$ws.Range("E10").Value = "'  -3.42%  "
$ws.Range("E11").Value = "'  -0.88%  "
$ws.Range("D12").Value = "'1.776.03"
$ws.Range("E12").Value = "'  -3.70%  "
$ws.Range("D13").Value = "'1.565.15"
$ws.Range("E13").Value = "'  -3.35%  "
$ws.Range("E14").Value = "'  -3.34%  "
$ws.Range("E15").Value = "'  -3.10%  "
$ws.Range("D16").Value = "'25.261.36"
$ws.Range("E16").Value = "'  -2.36%  "
$ws.Range("B17").Value = "'Litecoin"
$ws.Range("C17").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'59.19"
$ws.Range("E17").Value = "'  -3.56%  "
$ws.Range("B18").Value = "'ShibaInu"
$ws.Range("C18").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.0₃0711"
$ws.Range("E18").Value = "'  -3.46%  "
$ws.Range("E19").Value = "'  -0.24%  "
$ws.Range("D20").Value = "'186.44"
$ws.Range("E20").Value = "'  -2.56%  "
$ws.Range("E21").Value = "'  -2.89%  "
$ws.Range("E22").Value = "'  -2.39%  "
$ws.Range("E23").Value = "'  -2.64%  "
$ws.Range("E24").Value = "'  -2.84%  "
$ws.Range("D25").Value = "'141.02"
$ws.Range("E25").Value = "'  -1.64%  "
$ws.Range("E26").Value = "'  -0.21%  "
$ws.Range("E27").Value = "'  -2.49%  "
$ws.Range("E28").Value = "'  -1.60%  "
$ws.Range("D29").Value = "'6.38"
$ws.Range("E29").Value = "'  -4.32%  "
$ws.Range("E30").Value = "'  -7.15%  "
$ws.Range("D31").Value = "'0.0465"
$ws.Range("E31").Value = "'  -2.44%  "
$ws.Range("E32").Value = "'  -2.15%  "
$ws.Range("E33").Value = "'  -3.87%  "
$ws.Range("E34").Value = "'  -0.31%  "
$ws.Range("E35").Value = "'  -4.29%  "
$ws.Range("D36").Value = "'1.087.00"
$ws.Range("E36").Value = "'  -3.18%  "
$ws.Range("D38").Value = "'2.34"
$ws.Range("E38").Value = "'  -1.09%  "
$ws.Range("E39").Value = "'  -3.12%  "
$ws.Range("E40").Value = "'  -3.13%  "
$ws.Range("D41").Value = "'0.770"
$ws.Range("E41").Value = "'  -7.79%  "
$ws.Range("D42").Value = "'0.795"
$ws.Range("E42").Value = "'  +6.38%  "
$ws.Range("D43").Value = "'93.21"
$ws.Range("E43").Value = "'  -5.02%  "
$ws.Range("D44").Value = "'5.09"
$ws.Range("E44").Value = "'  +1.08%  "
$ws.Range("D45").Value = "'1.692.15"
$ws.Range("E45").Value = "'  -3.57%  "
$ws.Range("E46").Value = "'  -0.74%  "
$ws.Range("E47").Value = "'  -1.48%  "
$ws.Range("D48").Value = "'52.47"
$ws.Range("E48").Value = "'  -2.80%  "
$ws.Range("E49").Value = "'  -3.15%  "
$ws.Range("E50").Value = "'  -0.10%  "
$ws.Range("E51").Value = "'  -2.08%  "
